# Week 1 submission form update for BetfairDatathonPuntingFormDataSpecifications:
# - "Betfair Price" field description becomes "Betfair Starting Price (not available before the race)"
# - "Official starting price" field description becomes "Official starting price (or Pre-Post price if before the race)"
# (Order matters for shared-string table layout: Betfair entry is written first.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B47").Value = "Betfair Starting Price (not available before the race)"
$ws.Range("B45").Value = "Official starting price (or Pre-Post price if before the race)"

# Reflect the author's last on-screen position/selection (scrolled to row 38, B46 selected).
[void]$ws.Range("B46").Select()
